# ============================================================
# Insert a new "2022-Q4" worksheet right after "总计" and before "2022-Q3",
# populate it with the fund-holding table, and update the summary sheet.
# ============================================================
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)
$q3sheet = $wb.Worksheets.Item(2)

# --- Create the new "2022-Q4" sheet positioned before the "2022-Q3" sheet ---
$ws = $wb.Worksheets.Add($q3sheet)
$ws.Name = "2022-Q4"

# --- Header row (row 1) ---
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# --- Data rows (row 2 onward) ---
$ws.Range("A2").Value = 0
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "014562"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "易方达品质动能三年持有混合A"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "83.06"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "73.76"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "1.60"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "1.3290"
$ws.Range("H2").Value = 10
$ws.Range("A3").Value = 1
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "519732"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "交银定期支付双息平衡混合"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.61"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "69.68"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2.80"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "1.1091"
$ws.Range("H3").Value = 4
$ws.Range("A4").Value = 2
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "110029"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "易方达科讯混合"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "31.95"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "85.15"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "3.36"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "1.0735"
$ws.Range("H4").Value = 7
$ws.Range("A5").Value = 3
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "213001"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "宝盈鸿利收益灵活配置混合A"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "14.47"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "91.18"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "5.58"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "0.8074"
$ws.Range("H5").Value = 2
$ws.Range("A6").Value = 4
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "001076"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "易方达改革红利混合"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "25.93"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "93.28"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "2.79"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "0.7234"
$ws.Range("H6").Value = 8
$ws.Range("A7").Value = 5
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = "006533"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "易方达科融混合"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "18.30"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "88.33"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "3.45"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "0.6314"
$ws.Range("H7").Value = 8
$ws.Range("A8").Value = 6
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = "162202"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = "泰达宏利周期混合"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.93"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "89.10"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "2.76"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "0.2189"
$ws.Range("H8").Value = 8
$ws.Range("A9").Value = 7
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "014563"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "易方达品质动能三年持有混合C"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "12.59"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "73.76"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "1.60"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "0.2014"
$ws.Range("H9").Value = 10
$ws.Range("A10").Value = 8
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "003601"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "申万菱信安鑫精选混合A"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.71"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "36.51"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "4.82"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "0.1788"
$ws.Range("H10").Value = 2
$ws.Range("A11").Value = 9
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "005876"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "易方达鑫转增利混合A"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "29.60"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "1.10"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "0.0578"
$ws.Range("H11").Value = 8
$ws.Range("A12").Value = 10
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "005933"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "新疆前海联合先进制造灵活配置混合A"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.90"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "90.73"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "4.83"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.0435"
$ws.Range("H12").Value = 7
$ws.Range("A13").Value = 11
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "013877"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "财通资管新能源汽车混合C"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.74"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "94.54"
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "5.67"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "0.0420"
$ws.Range("H13").Value = 6
$ws.Range("A14").Value = 12
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "005877"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "易方达鑫转增利混合C"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.24"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "29.60"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "1.10"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "0.0356"
$ws.Range("H14").Value = 8
$ws.Range("A15").Value = 13
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "007581"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "宝盈鸿利收益灵活配置混合C"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.62"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "91.18"
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "5.58"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "0.0346"
$ws.Range("H15").Value = 2
$ws.Range("A16").Value = 14
$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = "004265"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "金鹰民丰回报定期开放混合"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.49"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "29.44"
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "0.73"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "0.0328"
$ws.Range("H16").Value = 7
$ws.Range("A17").Value = 15
$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = "011523"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = "新疆前海联合产业趋势混合A"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.60"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "83.30"
$ws.Range("F17").NumberFormat = "@"
$ws.Range("F17").Value = "4.75"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "0.0285"
$ws.Range("H17").Value = 4
$ws.Range("A18").Value = 16
$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = "011524"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = "新疆前海联合产业趋势混合C"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "83.30"
$ws.Range("F18").NumberFormat = "@"
$ws.Range("F18").Value = "4.75"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "0.0171"
$ws.Range("H18").Value = 4
$ws.Range("A19").Value = 17
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "013876"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "财通资管新能源汽车混合A"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.14"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "94.54"
$ws.Range("F19").NumberFormat = "@"
$ws.Range("F19").Value = "5.67"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "0.0079"
$ws.Range("H19").Value = 6
$ws.Range("A20").Value = 18
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "005934"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "新疆前海联合先进制造灵活配置混合C"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.09"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "90.73"
$ws.Range("F20").NumberFormat = "@"
$ws.Range("F20").Value = "4.83"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "0.0043"
$ws.Range("H20").Value = 7
$ws.Range("A21").Value = 19
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "009054"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "圆信永丰沣泰混合"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.23"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "31.90"
$ws.Range("F21").NumberFormat = "@"
$ws.Range("F21").Value = "1.20"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "0.0028"
$ws.Range("H21").Value = 8
$ws.Range("A22").Value = 20
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "005005"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "中金金泽量化精选混合A"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.09"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "69.88"
$ws.Range("F22").NumberFormat = "@"
$ws.Range("F22").Value = "2.81"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "0.0025"
$ws.Range("H22").Value = 5
$ws.Range("A23").Value = 21
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "005006"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = "中金金泽量化精选混合C"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.08"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "69.88"
$ws.Range("F23").NumberFormat = "@"
$ws.Range("F23").Value = "2.81"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "0.0022"
$ws.Range("H23").Value = 5
$ws.Range("A24").Value = 22
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "003602"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "申万菱信安鑫精选混合C"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "36.51"
$ws.Range("F24").NumberFormat = "@"
$ws.Range("F24").Value = "4.82"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "0.0019"
$ws.Range("H24").Value = 2

# --- Formatting: header row + column A (index) get bold font, thin border, center/top alignment ---
$hdr = $ws.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160

$idxCol = $ws.Range("A2:A24")
$idxCol.Font.Bold = $true
$idxCol.Borders.LineStyle = 1
$idxCol.HorizontalAlignment = -4108
$idxCol.VerticalAlignment = -4160

$ws.Range("A1").Select()

# ============================================================
# Update the "总计" (summary) sheet: insert the new 2022-Q4 row at the
# top of the data, and shift the rest of the quarters down by one.
# ============================================================
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 23
$summary.Range("D2").Value = 6.59
$summary.Range("B3").Value = "2022-Q3"
$summary.Range("C3").Value = 98
$summary.Range("D3").Value = 20.68
$summary.Range("B4").Value = "2022-Q2"
$summary.Range("C4").Value = 48
$summary.Range("D4").Value = 14.22
$summary.Range("B5").Value = "2022-Q1"
$summary.Range("C5").Value = 16
$summary.Range("D5").Value = 12.92
$summary.Range("B6").Value = "2021-Q4"
$summary.Range("C6").Value = 37
$summary.Range("D6").Value = 17.71
$summary.Range("B7").Value = "2021-Q3"
$summary.Range("C7").Value = 199
$summary.Range("D7").Value = 68.44
$summary.Range("B8").Value = "2021-Q2"
$summary.Range("C8").Value = 73
$summary.Range("D8").Value = 11.45
$summary.Range("B9").Value = "2021-Q1"
$summary.Range("C9").Value = 36
$summary.Range("D9").Value = 16.23
$summary.Range("B10").Value = "2020-Q4"
$summary.Range("C10").Value = 4
$summary.Range("D10").Value = 0.23

$summary.Range("A10").Value = 8

Write-Host "Edit complete"
